# Insert two new weekly price rows for "Cebolla" (1a/2a guarda) at the top
# of the existing data block (old row 264), pushing the rest of the table
# down by two rows (old 264..342 -> new 266..344).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 264-342 down by inserting two blank rows at 264.
$ws.Rows("264:265").Insert()

# Row 264: "1a (guarda)"
$ws.Range("A264").Value = 7
$ws.Range("B264").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C264").Value = "Ñuble"
$ws.Range("D264").Value = 44463
$ws.Range("E264").Value = 16
$ws.Range("F264").Value = 100112004
$ws.Range("G264").Value = "Cebolla"
$ws.Range("H264").Value = "Sin especificar"
$ws.Range("I264").Value = "1a (guarda)"
$ws.Range("J264").Value = 300
$ws.Range("K264").Value = 3000
$ws.Range("L264").Value = 3200
$ws.Range("M264").Value = 3100
$ws.Range("N264").Value = "$/malla 15 kilos"
$ws.Range("O264").Value = "Región del Maule"
$ws.Range("P264").Value = 207
$ws.Range("Q264").Value = 15
$ws.Range("R264").Value = "Hortaliza"

# Row 265: "2a (guarda)"
$ws.Range("A265").Value = 7
$ws.Range("B265").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C265").Value = "Ñuble"
$ws.Range("D265").Value = 44463
$ws.Range("E265").Value = 16
$ws.Range("F265").Value = 100112004
$ws.Range("G265").Value = "Cebolla"
$ws.Range("H265").Value = "Sin especificar"
$ws.Range("I265").Value = "2a (guarda)"
$ws.Range("J265").Value = 120
$ws.Range("K265").Value = 2800
$ws.Range("L265").Value = 2800
$ws.Range("M265").Value = 2800
$ws.Range("N265").Value = "$/malla 15 kilos"
$ws.Range("O265").Value = "Región del Maule"
$ws.Range("P265").Value = 187
$ws.Range("Q265").Value = 15
$ws.Range("R265").Value = "Hortaliza"

# Copy the date-column style (s="2") from the row above onto the new D cells,
# matching the rest of the "Fecha" column formatting.
$ws.Range("D263").Copy()
$ws.Range("D264:D265").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
